# Add 2022-Q4 data:
#  - new detail sheet "2022-Q4" inserted right after "总计" (and before "2022-Q3")
#  - the old "2022-Q3" detail sheet is preserved, shifted one position to the right
#  - "2022-Q2" / "2022-Q1" detail sheets are untouched, just shifted
#  - the "总计" summary sheet gets a new row for 2022-Q4 and the older rows shift down

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the current "2022-Q3" sheet right after itself -----------
# The duplicate keeps the old Q3 figures and will become the sheet named
# "2022-Q3" again; the original physical sheet becomes "2022-Q4" and gets the
# new quarter's numbers written into it.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $q3)

$q3.Name = "2022-Q4"

$q3Dup = $wb.Worksheets.Item("2022-Q3 (2)")
$q3Dup.Name = "2022-Q3"

# --- 2. Write the new 2022-Q4 figures into the (renamed) sheet -------------
# These "numeric-looking" figures are stored as TEXT in the source workbook
# (e.g. "4.20" keeps its trailing zero), so force text format before writing
# them, otherwise Excel auto-converts the strings to real numbers.
$q4 = $wb.Worksheets.Item("2022-Q4")

$q4.Range("E2:G2").NumberFormat = "@"
$q4.Range("E2").Value = "71.26"
$q4.Range("F2").Value = "4.81"
$q4.Range("G2").Value = "0.0034"
$q4.Range("H2").Value = 2

$q4.Range("D3:G3").NumberFormat = "@"
$q4.Range("D3").Value = "0.04"
$q4.Range("E3").Value = "71.26"
$q4.Range("F3").Value = "4.81"
$q4.Range("G3").Value = "0.0019"
$q4.Range("H3").Value = 2

# --- 3. Update the "总计" summary sheet -------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Shift existing rows 2-4 down to 3-5 (copy keeps number formats / styles)
$summary.Range("A4:D4").Copy($summary.Range("A5:D5"))
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# New row for 2022-Q4
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.01

# Fix up the running index column
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
